$wb = $excel.ActiveWorkbook
$aboutWs = $wb.Worksheets.Item("About")

# Remove the hyperlink on B6 (revert to plain text cell, keep existing formatting)
$aboutWs.Hyperlinks.Delete()

# Restore the old URL text value for B6
$aboutWs.Range("B6").Value = "http://yosemite.epa.gov/EE%5Cepa%5Ceed.nsf/webpages/MortalityRiskValuation.html#whatvalue"

# Restore the workbook calculation settings (disable iterative calculation)
$excel.Iteration = $false
